$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 changes
$ws.Range("K2").Value = 66.5
$ws.Range("N2").Value = 51.15965480231979

# Row 3 changes
$ws.Range("K3").Value = 65.3
$ws.Range("N3").Value = 51.15965480231979

# Row 4 changes
$ws.Range("D4").Value = 4255
$ws.Range("E4").Value = 72.8
$ws.Range("F4").Value = 0.87
$ws.Range("I4").Value = 73
$ws.Range("J4").Value = 76
$ws.Range("K4").Value = 59.5
$ws.Range("M4").Value = "⛔ 관망하십시오."
$ws.Range("N4").Value = 51.15965480231979
